$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 117 (shifts existing rows 117-171 down to 118-172),
# then populate the new row 117 with a new price observation.
$ws.Rows.Item(117).Insert()

$ws.Cells.Item(117, 1).Value = 7
$ws.Cells.Item(117, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(117, 3).Value = "Ñuble"
$ws.Cells.Item(117, 4).Value = "10/20/2021"
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 6).Value = 100112008
$ws.Cells.Item(117, 7).Value = "Coliflor"
$ws.Cells.Item(117, 8).Value = "Sin especificar"
$ws.Cells.Item(117, 9).Value = "Primera"
$ws.Cells.Item(117, 10).Value = 180
$ws.Cells.Item(117, 11).Value = 700
$ws.Cells.Item(117, 12).Value = 800
$ws.Cells.Item(117, 13).Value = 750
$ws.Cells.Item(117, 14).Value = "$/unidad"
$ws.Cells.Item(117, 15).Value = "Región Metropolitana"
$ws.Cells.Item(117, 16).Value = 750
$ws.Cells.Item(117, 17).Value = 1
$ws.Cells.Item(117, 18).Value = "Hortaliza"
